# The commit changes cell B11 on the "Rules" sheet from the text "R40"
# to the text "1" (still a text/shared-string value, not a number — the
# OOXML keeps t="s" pointing at a (new) shared-string entry and leaves the
# cell's existing style (s="23") untouched).
#
# A plain `Range.Value = "1"` would make Excel's type-inference treat the
# digit string as a genuine number (t omitted, <v>1</v>), and prefixing
# with an apostrophe to force text (`'1`) flips the cell onto a brand new
# "quote prefixed" style variant instead of keeping the original style id.
# To reproduce the workbook's actual end state (text "1", original style
# untouched) we: stash the cell's current formatting on a scratch cell far
# outside the used range, assign the text value (which bumps B11 onto the
# auto-generated quote-prefix style), then paste the stashed formatting
# back over B11 so it reports its original style again. The scratch row is
# then deleted so no trace/residue is left in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target  = $ws.Cells.Item(11, 2)     # B11
$scratch = $ws.Cells.Item(100, 26)   # Z100 - far outside the used range

# 1) Remember B11's current formatting.
$target.Copy()
$scratch.PasteSpecial(-4122)   # xlPasteFormats

# 2) Write the new text value. Leading apostrophe forces it to be stored
#    as text ("1"), matching the shared-string t="s" cell in the target.
$target.Value = "'1"

# 3) Restore B11's original formatting/style (removes the quote-prefix
#    style Excel auto-assigned in step 2) without touching its value.
$scratch.Copy()
$target.PasteSpecial(-4122)    # xlPasteFormats

# 4) Clean up the scratch cell/row so it leaves no trace in the sheet.
$scratch.EntireRow.Delete()
$excel.CutCopyMode = $false
